# Add the newly finished book "Why We Sleep" to the "Completed" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# New row goes right after the last existing row of data (row 25 -> row 26).
$row = 26

$ws.Cells.Item($row, 1).Value = "Why We Sleep"
$ws.Cells.Item($row, 2).Value = "Matthew Walker"
$ws.Cells.Item($row, 3).Value = [DateTime]"2020-02-15"
$ws.Cells.Item($row, 4).Value = [DateTime]"2020-02-17"
$ws.Cells.Item($row, 5).Value = "sleep;health;science;sleep deprivation;disease;wellness"
$ws.Cells.Item($row, 6).Value = "Audio"
$ws.Cells.Item($row, 7).Value = "14 Hours 0 Mins"

# Match the date formatting already used in column C/D (reuse the existing
# style rather than minting a new number format) by copying formats down
# from existing date cells (row 25 has no Finish Date, so use row 24 for D).
$ws.Range("C25").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("D24").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$excel.CutCopyMode = 0
